# Apply "Thay đổi dữ liệu" changes to 06.DanhSachChucNang.xlsx
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 8: "Phản hồi" -> Hoàn thành 100%, Phân công "Tú, Nhi"
$ws.Range("E8").Value = 1
$ws.Range("F8").Value = "Tú, Nhi"

# Row 6: "Tìm kiếm chuyến xe" -> Hoàn thành 70%, Phân công "Huy, Kiều"
$ws.Range("E6").Value = 0.7
$ws.Range("F6").Value = "Huy, Kiều"

# Row 34: "Thêm chuyến" -> Hoàn thành 90%, Phân công "Lê, Kiều"
$ws.Range("E34").Value = 0.9
$ws.Range("F34").Value = "Lê, Kiều"

# Move the view / selection like the author did while editing
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 22
$ws.Range("F35").Select()
